$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 16-18: field/value pairs appended below the existing data
$ws.Range("A16").Value = "Invalid Date Error"
$ws.Range("B16").Value = "date must be a ``date`` type, but the final value was: ``Invalid Date``."

$ws.Range("A17").Value = "Invalid Visit Reason"
$ws.Range("B17").Value = "Field is required"

$ws.Range("A18").Value = "Invalid Birthdate"
$ws.Range("B18").Value = "birthdate must be a ``date`` type, but the final value was: ``Invalid Date``."

# Update the active selection to match the post-edit state
$ws.Range("D20").Select()
